# Handback status report regeneration: refresh the "Correspond Handback
# DateTime" / "Latest HO Xliff Generate Date" timestamps for the
# b99a1b6e-0378-487d-8ebd-43927b15cb35.md row on each sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-29 06:46:28"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-29 06:46:24"
$wsZhCn.Range("K3").Value = "2016-08-29 06:46:42"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-29 06:46:28"
$wsDeDe.Range("K3").Value = "2016-08-29 06:46:49"
